$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-19 with new contact data (Canadian Flatiron roster)
# Row 2
$ws.Range("A2").Value = "Paul Newman"
$ws.Range("B2").Value = "Paul"
$ws.Range("C2").Value = "Newman"
$ws.Range("D2").Value = "https://linkedin.com/in/paul-newman-00a9a830"
$ws.Range("E2").Value = "Bow Transit Connectors"
$ws.Range("F2").Value = "Construction Manager"
$ws.Range("G2").Value = "Calgary, Alberta, Canada"
$ws.Range("H2").Value = "+1 403-701-8000 , +1 403-589-9088 , +1 403-620-7239 , +1 410-340-7346 , +1 418-868-2428 , +1 403-640-3521"
$ws.Range("I2").Value = "paul.newman@snclavalin.com , paul@pave-it.com , pnewman@bellsouth.net"

# Row 3
$ws.Range("A3").Value = "Neil Lynchehaun"
$ws.Range("B3").Value = "Neil"
$ws.Range("C3").Value = "Lynchehaun"
$ws.Range("D3").Value = "https://linkedin.com/in/neil-lynchehaun-23031ba4"
$ws.Range("E3").Value = "Flatiron Construction"
$ws.Range("F3").Value = "Vice President, Special Projects"
$ws.Range("G3").Value = "Calgary, Alberta, Canada"

# Row 4
$ws.Range("A4").Value = "Ken Tanner"
$ws.Range("B4").Value = "Ken"
$ws.Range("C4").Value = "Tanner"
$ws.Range("D4").Value = "https://linkedin.com/in/ken-tanner-9899939"
$ws.Range("E4").Value = "Flatiron Construction"
$ws.Range("F4").Value = "Vice President Operations"
$ws.Range("G4").Value = "Canada"
$ws.Range("H4").Value = "+1 604-798-8491"
$ws.Range("I4").Value = "ktanner@flatironcorp.com"

# Row 5
$ws.Range("A5").Value = "Octavio Flores"
$ws.Range("B5").Value = "Octavio"
$ws.Range("C5").Value = "Flores"
$ws.Range("D5").Value = "https://linkedin.com/in/octavio-flores-37b566a7"
$ws.Range("E5").Value = "Flatiron Construction"
$ws.Range("F5").Value = "Area Operations Manager"
$ws.Range("G5").Value = "British Columbia, Canada"

# Row 6
$ws.Range("A6").Value = "Jarred Gumbleton"
$ws.Range("B6").Value = "Jarred"
$ws.Range("C6").Value = "Gumbleton"
$ws.Range("D6").Value = "https://linkedin.com/in/jarredgumbleton"
$ws.Range("E6").Value = "Flatiron Construction"
$ws.Range("F6").Value = "Project Manager"
$ws.Range("G6").Value = "Richmond, BC, Canada"
$ws.Range("H6").Value = "+1 250-272-6645"
$ws.Range("I6").Value = "jgumbleton@flatironcorp.com"

# Row 7
$ws.Range("A7").Value = "Alicia Lopez"
$ws.Range("B7").Value = "Alicia"
$ws.Range("C7").Value = "Lopez"
$ws.Range("D7").Value = "https://linkedin.com/in/alicialopezcrespo"
$ws.Range("E7").Value = "Flatiron Construction"
$ws.Range("F7").Value = "VP, Innovation & Design"
$ws.Range("G7").Value = "USA and Canada "
$ws.Range("H7").Value = "DNC , DNC , +1 719-994-0855 , +1 720-232-3886 , DNC"
$ws.Range("I7").Value = "acrespo@flatironcorp.com"

# Row 8
$ws.Range("A8").Value = "Mike M"
$ws.Range("B8").Value = "Mike"
$ws.Range("C8").Value = "M"
$ws.Range("D8").Value = "https://linkedin.com/in/mike-meacher-031ba619"
$ws.Range("E8").Value = "Flatiron Construction"
$ws.Range("F8").Value = "Quality Manager"
$ws.Range("G8").Value = "British Columbia, Canada"
$ws.Range("H8").Value = "+1 604-363-4558"
$ws.Range("I8").Value = "mike.meacher@lafarge-na.com , mmeacher@fwsgroup.com , mmeacher@flatironcorp.com"

# Row 9
$ws.Range("A9").Value = "Donald Dow"
$ws.Range("B9").Value = "Donald"
$ws.Range("C9").Value = "Dow"
$ws.Range("D9").Value = "https://linkedin.com/in/donald-dow-7518ab192"
$ws.Range("E9").Value = "AFDE Partnership"
$ws.Range("F9").Value = "Safety Manager"
$ws.Range("G9").Value = "Fort St John, British Columbia, Canada"

# Row 10
$ws.Range("A10").Value = "Jan Kyrstein"
$ws.Range("B10").Value = "Jan"
$ws.Range("C10").Value = "Kyrstein"
$ws.Range("D10").Value = "https://linkedin.com/in/jan-kyrstein-39007211"
$ws.Range("E10").Value = "Flatiron Construction"
$ws.Range("F10").Value = "Construction Manager/Deputy Area Manager "
$ws.Range("G10").Value = "Fort St John, British Columbia, Canada"
$ws.Range("H10").Value = "+1 306-318-0214 , +1 604-563-5197"
$ws.Range("I10").Value = "jan.kyrstein@pinnaclepellet.com , jan.kyrstein@bhpbilliton.com , kyrstein@technologist.com"

# Row 11
$ws.Range("A11").Value = "Dana Driver"
$ws.Range("B11").Value = "Dana"
$ws.Range("C11").Value = "Driver"
$ws.Range("D11").Value = "https://linkedin.com/in/dana-driver-crsp-9050584a"
$ws.Range("E11").Value = "Flatiron Construction"
$ws.Range("F11").Value = "District Safety Manager"
$ws.Range("G11").Value = "Keeyask"
$ws.Range("H11").Value = "DNC , +1 587-434-3086 , +1 403-280-2126"
$ws.Range("I11").Value = "ddriver@flatironcorp.com"

# Row 12
$ws.Range("A12").Value = "Mark Neis"
$ws.Range("B12").Value = "Mark"
$ws.Range("C12").Value = "Neis"
$ws.Range("D12").Value = "https://linkedin.com/in/mark-neis-27462812"
$ws.Range("E12").Value = "Flatiron Construction"
$ws.Range("F12").Value = "Deputy Project Director"
$ws.Range("G12").Value = "Fort St John, British Columbia, Canada"
$ws.Range("H12").Value = "+63 926 022 0867 , +63 942 908 8629"
$ws.Range("I12").Value = "100267.2650@compuserve.com , mark.neis@cbi.com , mneis@flatironcorp.com"

# Row 13
$ws.Range("A13").Value = "Joel Jacques"
$ws.Range("B13").Value = "Joel"
$ws.Range("C13").Value = "Jacques"
$ws.Range("D13").Value = "https://linkedin.com/in/joel-jacques-p-eng-7072a655"
$ws.Range("E13").Value = "Coast Valley Contracting Ltd"
$ws.Range("F13").Value = "President"
$ws.Range("G13").Value = "Squamish, British Columbia, Canada"
$ws.Range("H13").Value = "+1 604-849-1017"
$ws.Range("I13").Value = "joel.jacques@coastvalley.ca"

# Row 14
$ws.Range("A14").Value = "Fidel Velarde"
$ws.Range("B14").Value = "Fidel"
$ws.Range("C14").Value = "Velarde"
$ws.Range("D14").Value = "https://linkedin.com/in/fidel-velarde-054366164"
$ws.Range("E14").Value = "Flatiron Construction"
$ws.Range("F14").Value = "Quality Manager"
$ws.Range("G14").Value = "Richmond, British Columbia, Canada"

# Row 15
$ws.Range("A15").Value = "Jenn Hirschman"
$ws.Range("B15").Value = "Jenn"
$ws.Range("C15").Value = "Hirschman"
$ws.Range("D15").Value = "https://linkedin.com/in/jenn-hirschman-7b55a1167"
$ws.Range("E15").Value = "Flatiron Construction"
$ws.Range("F15").Value = "Health & Safety Manager III"
$ws.Range("G15").Value = "British Columbia, Canada"

# Row 16
$ws.Range("A16").Value = "Pete Walton"
$ws.Range("B16").Value = "Pete"
$ws.Range("C16").Value = "Walton"
$ws.Range("D16").Value = "https://linkedin.com/in/pete-walton"
$ws.Range("E16").Value = "Flatiron Construction"
$ws.Range("F16").Value = "Canadian Division Safety Director"
$ws.Range("G16").Value = "Alberta Construction Safety Association"
$ws.Range("H16").Value = "+1 970-290-8000 , +1 252-229-5369 , +1 303-579-7895 , +1 719-351-5126 , +1 224-578-1158 , +1 210-618-2993 , +1 956-467-6750"
$ws.Range("I16").Value = "plwalton@verizon.net , irmawalton@webtv.net"

# Row 17
$ws.Range("A17").Value = "Frank Mydlinski"
$ws.Range("B17").Value = "Frank"
$ws.Range("C17").Value = "Mydlinski"
$ws.Range("D17").Value = "https://linkedin.com/in/frank-mydlinski-0566a184"
$ws.Range("E17").Value = "Flatiron Construction"
$ws.Range("F17").Value = "Senior Quality Manager"
$ws.Range("G17").Value = "Chilliwack, British Columbia, Canada"
$ws.Range("H17").Value = ""
$ws.Range("I17").Value = ""

# Row 18
$ws.Range("A18").Value = "Husted Janet"
$ws.Range("B18").Value = "Husted"
$ws.Range("C18").Value = "Janet"
$ws.Range("D18").Value = "https://linkedin.com/in/husted-janet-62b1b8116"
$ws.Range("E18").Value = "Flatiron Construction"
$ws.Range("F18").Value = "Health Safety Manager"
$ws.Range("G18").Value = "Merritt, British Columbia, Canada"
$ws.Range("H18").Value = ""
$ws.Range("I18").Value = ""

# Row 19
$ws.Range("A19").Value = "Hodge Garry"
$ws.Range("B19").Value = "Hodge"
$ws.Range("C19").Value = "Garry"
$ws.Range("D19").Value = "https://linkedin.com/in/hodge-garry-a4a3a59b"
$ws.Range("E19").Value = "Flatiron Construction"
$ws.Range("F19").Value = "PROJECT SAFETY MANAGER"
$ws.Range("G19").Value = "Greater Edmonton Metropolitan Area"

# Remove rows 20-22 (no longer present in updated roster)
$ws.Rows.Item(20).EntireRow.Delete()
$ws.Rows.Item(20).EntireRow.Delete()
$ws.Rows.Item(20).EntireRow.Delete()
